$wb = $excel.ActiveWorkbook

# Update the workbook-wide revision pointer is left to the engine; we only
# need to add the new "Menu" worksheet after the last existing sheet ("156").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Menu"

# --- Column widths (approximate match to the target bestFit widths) ---
$ws.Columns.Item(3).ColumnWidth = 11.053385416666666
$ws.Columns.Item(4).ColumnWidth = 16.721354166666668
$ws.Columns.Item(5).ColumnWidth = 18.608072916666668
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668

# --- Cell values: written in the exact order the strings were first used ---
# (this preserves the shared-strings table order from the source edit)
$ws.Range("A2").Value = "Run"
$ws.Range("G5").Value = ">"
$ws.Range("B2").Value = ">login"
$ws.Range("C2").Value = ">Admin"
$ws.Range("E2").Value = ">Show product"
$ws.Range("F2").Value = ">Show all"
$ws.Range("F3").Value = ">Show Limited"
$ws.Range("F4").Value = ">Show unlimited"
$ws.Range("F5").Value = ">sort by"
$ws.Range("E7").Value = ">Create note"
$ws.Range("F7").Value = ">receive note"
$ws.Range("F8").Value = ">delivery note"
$ws.Range("F10").Value = ">Remove product"
$ws.Range("E10").Value = ">Find Product"
$ws.Range("D14").Value = ">Business Manager"
$ws.Range("D2").Value = ">Product Manager"
$ws.Range("E14").Value = ">Show note list"
$ws.Range("E18").Value = ">calculate"
$ws.Range("F6").Value = ">back to previous"
$ws.Range("E13").Value = ">Back to previous"
$ws.Range("F14").Value = ">show receive note"
$ws.Range("F15").Value = ">show delivery note"
$ws.Range("E20").Value = ">add money"
$ws.Range("D23").Value = ">User manager"
$ws.Range("F16").Value = ">find note"
$ws.Range("E23").Value = ">show user list"
$ws.Range("F23").Value = ">show all"
$ws.Range("F24").Value = ">show accountant"
$ws.Range("F25").Value = ">show store keeper"
$ws.Range("F26").Value = ">show sale staff "
$ws.Range("E28").Value = ">create new user "
$ws.Range("F28").Value = ">accountant"
$ws.Range("F29").Value = ">storekeeper"
$ws.Range("F30").Value = ">salestaff"
$ws.Range("E33").Value = ">remove user"
$ws.Range("E32").Value = ">find user"
$ws.Range("C37").Value = ">Accountant"
$ws.Range("C50").Value = ">Storekeeper"
$ws.Range("C62").Value = ">SaleStaff"
$ws.Range("D35").Value = ">Log out"
$ws.Range("B74").Value = ">exit"

# --- Remaining cells that reuse already-interned strings ---
$ws.Range("F9").Value = ">back to previous"
$ws.Range("F11").Value = ">back to previous"
$ws.Range("E12").Value = ">Remove product"
$ws.Range("F17").Value = ">back to previous"
$ws.Range("E21").Value = ">Back to previous"
$ws.Range("F27").Value = ">back to previous"
$ws.Range("F31").Value = ">back to previous"
$ws.Range("E34").Value = ">Back to previous"
$ws.Range("D37").Value = ">Show product"
$ws.Range("E37").Value = ">Show all"
$ws.Range("E38").Value = ">Show Limited"
$ws.Range("E39").Value = ">Show unlimited"
$ws.Range("E40").Value = ">sort by"
$ws.Range("E41").Value = ">back to previous"
$ws.Range("D42").Value = ">Show note list"
$ws.Range("E42").Value = ">show receive note"
$ws.Range("E43").Value = ">show delivery note"
$ws.Range("E44").Value = ">find note"
$ws.Range("E45").Value = ">back to previous"
$ws.Range("D46").Value = ">calculate"
$ws.Range("D47").Value = ">Back to previous"
$ws.Range("D48").Value = ">Log out"
$ws.Range("D50").Value = ">Show product"
$ws.Range("E50").Value = ">Show all"
$ws.Range("E51").Value = ">Show Limited"
$ws.Range("E52").Value = ">Show unlimited"
$ws.Range("E53").Value = ">sort by"
$ws.Range("E54").Value = ">back to previous"
$ws.Range("D55").Value = ">Create note"
$ws.Range("E55").Value = ">receive note"
$ws.Range("E56").Value = ">delivery note"
$ws.Range("E57").Value = ">back to previous"
$ws.Range("D58").Value = ">Find Product"
$ws.Range("E58").Value = ">Remove product"
$ws.Range("E59").Value = ">back to previous"
$ws.Range("D60").Value = ">Log out"
$ws.Range("D62").Value = ">Show product"
$ws.Range("E62").Value = ">Show all"
$ws.Range("E63").Value = ">Show Limited"
$ws.Range("E64").Value = ">Show unlimited"
$ws.Range("E65").Value = ">sort by"
$ws.Range("E66").Value = ">back to previous"
$ws.Range("D67").Value = ">Create note"
$ws.Range("E67").Value = ">delivery note"
$ws.Range("D68").Value = ">show delivery note"
$ws.Range("D69").Value = ">Back to previous"
$ws.Range("D70").Value = ">Log out"
$ws.Range("C72").Value = ">Log out"

# --- Row heights for the two explicitly-sized rows near the bottom ---
$ws.Rows.Item(70).RowHeight = 15
$ws.Rows.Item(71).RowHeight = 15

# --- Final selection/view state ---
$ws.Range("E18").Select()
